# Update generated view-count figures (column F) on the "展览" and
# "全部类型" sheets to match the freshly generated data output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2084
$ws1.Range("F4").Value = 860
$ws1.Range("F5").Value = 1242
$ws1.Range("F6").Value = 359

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2084
$ws4.Range("F6").Value = 860
$ws4.Range("F7").Value = 1242
$ws4.Range("F8").Value = 359
